$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.081.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.299.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.60%  '
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.86'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.656.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.319.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.950.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.64%  '
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.18%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.93'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.41'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.011.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('E47').Value = '  -2.29%  '
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.523.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.54%  '
